$wb = $excel.ActiveWorkbook

# --- 1. Status text update -------------------------------------------------
# "Ready for handoff" -> "In Translation" everywhere it appears:
#   Overview!E2:F4 (zh-cn / de-de status columns on the summary sheet)
#   zh-cn!C2:C4    (Status column)
#   de-de!C2:C4    (Status column)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C4").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C4").Value = "In Translation"

# --- 2. Shrink the (now shorter) status columns -----------------------------
# The status text went from "Ready for handoff" (18 chars) to
# "In Translation" (15 chars), so the columns that display it are narrowed
# to fit the new content.
$wsOverview.Range("E1:F1").EntireColumn.ColumnWidth = 12.5
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = 12.5
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = 12.5
